# Product Backlog v1.xlsx - "atualizado status e backlog"
#
# The PBC sheet has a set of requirement rows (9-16) whose short "title"
# (column C) all used to be prefixed "Oficina - ...". They are renamed to
# use the "Serviços - ..." prefix instead (the long descriptions in column D
# are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PBC")
$ws.Activate()

# Order matches the order the strings were (re)introduced in the workbook's
# shared-string table so the new entries are appended the same way.
$ws.Range("C11").Value = "Serviços - Extrato"
$ws.Range("C12").Value = "Serviços - Telas Detalhe"
$ws.Range("C13").Value = "Serviços - Filtro"
$ws.Range("C14").Value = "Serviços - Cadastro Cliente"
$ws.Range("C15").Value = "Serviços - Botão Nova Aplicação"
$ws.Range("C16").Value = "Serviços - Tela Aplicação"
$ws.Range("C10").Value = "Serviços - Dashboard Mais Detalhes"
$ws.Range("C9").Value  = "Serviços - Dashboard Gráfico"

# Move the cursor / visible window to where the edits were made.
[void]$ws.Range("C12").Select()
